# This script reproduces a commit that inserted two new price-record rows
# (rows 1109 and 1110) into the "Fruta, Vega Modelo de Temuco - Plátano"
# daily price sheet, pushing all subsequent rows down by two positions
# (old row 1109 -> 1111, old row 1110 -> 1112, ..., old row 1201 -> 1203).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 1109, shifting existing data down.
$ws.Rows("1109:1110").Insert()

# --- Fill in new row 1109 ---
$ws.Cells.Item(1109, 1).Value = 10
$ws.Cells.Item(1109, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(1109, 3).Value = 'La Araucanía'
$ws.Cells.Item(1109, 4).Value = 45106
$ws.Cells.Item(1109, 5).Value = 9
$ws.Cells.Item(1109, 6).Value = 'Fruta'
$ws.Cells.Item(1109, 7).Value = 100108
$ws.Cells.Item(1109, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(1109, 9).Value = 100108006
$ws.Cells.Item(1109, 10).Value = 'Plátano'
$ws.Cells.Item(1109, 11).Value = 'Barraganete'
$ws.Cells.Item(1109, 12).Value = 'Primera'
$ws.Cells.Item(1109, 13).Value = 80
$ws.Cells.Item(1109, 14).Value = 36000
$ws.Cells.Item(1109, 15).Value = 36000
$ws.Cells.Item(1109, 16).Value = 36000
$ws.Cells.Item(1109, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(1109, 18).Value = 'Ecuador'
$ws.Cells.Item(1109, 19).Value = 1800
$ws.Cells.Item(1109, 20).Value = 20

# --- Fill in new row 1110 ---
$ws.Cells.Item(1110, 1).Value = 10
$ws.Cells.Item(1110, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(1110, 3).Value = 'La Araucanía'
$ws.Cells.Item(1110, 4).Value = 45106
$ws.Cells.Item(1110, 5).Value = 9
$ws.Cells.Item(1110, 6).Value = 'Fruta'
$ws.Cells.Item(1110, 7).Value = 100108
$ws.Cells.Item(1110, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(1110, 9).Value = 100108006
$ws.Cells.Item(1110, 10).Value = 'Plátano'
$ws.Cells.Item(1110, 11).Value = 'Sin especificar'
$ws.Cells.Item(1110, 12).Value = 'Pintón'
$ws.Cells.Item(1110, 13).Value = 2200
$ws.Cells.Item(1110, 14).Value = 16000
$ws.Cells.Item(1110, 15).Value = 17000
$ws.Cells.Item(1110, 16).Value = 16455
$ws.Cells.Item(1110, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(1110, 18).Value = 'Ecuador'
$ws.Cells.Item(1110, 19).Value = 823
$ws.Cells.Item(1110, 20).Value = 20

# Make sure the Date column (D) keeps the date/time number format used by
# the rest of the column (style index 2 in the original workbook).
$ws.Range("D1109:D1110").NumberFormat = $ws.Range("D1108").NumberFormat

Write-Output "Inserted rows 1109-1110; new dimension should be A1:T1203"
